# rebuild project Use Unirest to send request
#
# Adds a new data row to the "TZ" sheet (row 3) mirroring the existing
# request-field row, with a new license-plate value, and updates the
# active-cell selection on the "TZ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TZ")

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "АВ2991АЕ"

$ws.Range("F7").Select()
